# Auto-generated edit script: update Leve profit calc sheets
# per scheduled runner price refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 81
$ws.Range("H81").Value = 32928
$ws.Range("J81").Value = 32928
$ws.Range("L81").Value = 32928
$ws.Range("N81").Value = -34924

# Row 84
$ws.Range("H84").Value = 32928
$ws.Range("J84").Value = 32928
$ws.Range("L84").Value = 98784
$ws.Range("N84").Value = -108768

# Row 92
$ws.Range("H92").Value = 480.13333
$ws.Range("I92").Value = 469.3846
$ws.Range("K92").Value = 469.3846
$ws.Range("M92").Value = 778.6154

# Row 137
$ws.Range("H137").Value = 1222.7693
$ws.Range("I137").Value = 1127.1
$ws.Range("J137").Value = 1541.6666
$ws.Range("K137").Value = 3381.3
$ws.Range("L137").Value = 4624.9998
$ws.Range("M137").Value = -831.2999999999997
$ws.Range("N137").Value = -9724.9998

# Row 141
$ws.Range("H141").Value = 2725.8635
$ws.Range("I141").Value = 2566.4614
$ws.Range("K141").Value = 7699.3842
$ws.Range("M141").Value = -2519.3842

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5755.9673
$ws.Range("I32").Value = 4465.7534
$ws.Range("K32").Value = 4465.7534
$ws.Range("M32").Value = -4178.7534

# Row 61
$ws.Range("H61").Value = 4000.4546
$ws.Range("I61").Value = 4263.125
$ws.Range("K61").Value = 4263.125
$ws.Range("M61").Value = -4051.125

# Row 63
$ws.Range("H63").Value = 2842650
$ws.Range("I63").Value = 1915.1
$ws.Range("K63").Value = 1915.1
$ws.Range("M63").Value = -1229.1

# Row 66
$ws.Range("H66").Value = 2842650
$ws.Range("I66").Value = 1915.1
$ws.Range("K66").Value = 9575.5
$ws.Range("M66").Value = -6143.5

# Row 97
$ws.Range("H97").Value = 58824900
$ws.Range("I97").Value = 1008.8182
$ws.Range("K97").Value = 1008.8182
$ws.Range("M97").Value = -512.8182

# Row 136
$ws.Range("H136").Value = 4000.4546
$ws.Range("I136").Value = 4263.125
$ws.Range("K136").Value = 12789.375
$ws.Range("M136").Value = -10239.375

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 722.0625
$ws.Range("I22").Value = 679.36365
$ws.Range("J22").Value = 816
$ws.Range("K22").Value = 679.36365
$ws.Range("L22").Value = 816
$ws.Range("M22").Value = -506.36365
$ws.Range("N22").Value = -1162

# Row 35
$ws.Range("H35").Value = 24136.334
$ws.Range("J35").Value = 24136.334
$ws.Range("L35").Value = 24136.334
$ws.Range("N35").Value = -24756.334

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 17.142857
$ws.Range("I7").Value = 11.444445
$ws.Range("J7").Value = 27.4
$ws.Range("K7").Value = 11.444445
$ws.Range("L7").Value = 27.4
$ws.Range("M7").Value = 101.555555
$ws.Range("N7").Value = -253.4

# Row 16
$ws.Range("H16").Value = 2333.3333
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3074

# Row 31
$ws.Range("H31").Value = 4695.731
$ws.Range("I31").Value = 2521
$ws.Range("J31").Value = 7232.9165
$ws.Range("K31").Value = 2521
$ws.Range("L31").Value = 7232.9165
$ws.Range("M31").Value = -2226
$ws.Range("N31").Value = -7822.9165

# Row 34
$ws.Range("H34").Value = 4695.731
$ws.Range("I34").Value = 2521
$ws.Range("J34").Value = 7232.9165
$ws.Range("K34").Value = 2521
$ws.Range("L34").Value = 7232.9165
$ws.Range("M34").Value = -2319
$ws.Range("N34").Value = -7636.9165

# Row 51
$ws.Range("H51").Value = 34450
$ws.Range("J51").Value = 34450
$ws.Range("L51").Value = 34450
$ws.Range("N51").Value = -35922

# Row 61
$ws.Range("H61").Value = 34450
$ws.Range("J61").Value = 34450
$ws.Range("L61").Value = 34450
$ws.Range("N61").Value = -35146

# Row 99
$ws.Range("H99").Value = 3101.963
$ws.Range("I99").Value = 2685.5557
$ws.Range("J99").Value = 3934.7778
$ws.Range("K99").Value = 2685.5557
$ws.Range("L99").Value = 3934.7778
$ws.Range("M99").Value = -1187.5557
$ws.Range("N99").Value = -6930.7778

# Row 113
$ws.Range("H113").Value = 2333.3333
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -6840

# Row 126
$ws.Range("H126").Value = 3101.963
$ws.Range("I126").Value = 2685.5557
$ws.Range("J126").Value = 3934.7778
$ws.Range("K126").Value = 8056.6671
$ws.Range("L126").Value = 11804.3334
$ws.Range("M126").Value = -5586.6671
$ws.Range("N126").Value = -16744.3334

# Row 134
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 12522.556
$ws.Range("I68").Value = 1150
$ws.Range("J68").Value = 15771.857
$ws.Range("K68").Value = 3450
$ws.Range("L68").Value = 47315.571
$ws.Range("M68").Value = -2639
$ws.Range("N68").Value = -48937.571

# Row 71
$ws.Range("H71").Value = 12522.556
$ws.Range("I71").Value = 1150
$ws.Range("J71").Value = 15771.857
$ws.Range("K71").Value = 10350
$ws.Range("L71").Value = 141946.713
$ws.Range("M71").Value = -6294
$ws.Range("N71").Value = -150058.713

# Row 86
$ws.Range("H86").Value = 1450
$ws.Range("I86").Value = 662.5
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1987.5
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = -801.5
$ws.Range("N86").Value = -9872

# Row 89
$ws.Range("H89").Value = 1450
$ws.Range("I89").Value = 662.5
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 5962.5
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -34.5
$ws.Range("N89").Value = -34356

# Row 111
$ws.Range("H111").Value = 238.5
$ws.Range("I111").Value = 238.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 715.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 2351.5
$ws.Range("N111").ClearContents()

# Row 131
$ws.Range("H131").Value = 704.61
$ws.Range("J131").Value = 726.1398
$ws.Range("L131").Value = 2178.4194
$ws.Range("N131").Value = -12258.4194

# Row 139
$ws.Range("H139").Value = 2268.205
$ws.Range("I139").Value = 1428
$ws.Range("J139").Value = 3476
$ws.Range("K139").Value = 4284
$ws.Range("L139").Value = 10428
$ws.Range("M139").Value = 856
$ws.Range("N139").Value = -20708

# Row 141
$ws.Range("H141").Value = 3665
$ws.Range("I141").Value = 2272.5
$ws.Range("K141").Value = 6817.5
$ws.Range("M141").Value = -1637.5

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 20142.857
$ws.Range("J46").Value = 20166.666
$ws.Range("L46").Value = 20166.666
$ws.Range("N46").Value = -20478.666

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2814.074
$ws.Range("I40").Value = 2510.25
$ws.Range("J40").Value = 3682.1428
$ws.Range("K40").Value = 2510.25
$ws.Range("L40").Value = 3682.1428
$ws.Range("M40").Value = -2374.25
$ws.Range("N40").Value = -3954.1428

# Row 46
$ws.Range("H46").Value = 1742.2759
$ws.Range("I46").Value = 1774.5769
$ws.Range("J46").Value = 1462.3334
$ws.Range("K46").Value = 1774.5769
$ws.Range("L46").Value = 1462.3334
$ws.Range("M46").Value = -1586.5769
$ws.Range("N46").Value = -1838.3334

# Row 93
$ws.Range("H93").Value = 1192.7273
$ws.Range("I93").Value = 1192.7273
$ws.Range("K93").Value = 1192.7273
$ws.Range("M93").Value = 55.27269999999999

# Row 122
$ws.Range("H122").Value = 757716.6
$ws.Range("I122").Value = 936730.1
$ws.Range("K122").Value = 2810190.3
$ws.Range("M122").Value = -2807740.3

# Row 140
$ws.Range("H140").Value = 48429
$ws.Range("J140").Value = 48429
$ws.Range("L140").Value = 48429
$ws.Range("N140").Value = -58789

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2213.6667
$ws.Range("I81").Value = 1577.6
$ws.Range("J81").Value = 2668
$ws.Range("K81").Value = 3155.2
$ws.Range("L81").Value = 5336
$ws.Range("M81").Value = -2094.2
$ws.Range("N81").Value = -7458

# Row 84
$ws.Range("H84").Value = 2213.6667
$ws.Range("I84").Value = 1577.6
$ws.Range("J84").Value = 2668
$ws.Range("K84").Value = 15776
$ws.Range("L84").Value = 26680
$ws.Range("M84").Value = -10472
$ws.Range("N84").Value = -37288

# Row 132
$ws.Range("H132").Value = 1275.75
$ws.Range("I132").Value = 1052.931
$ws.Range("J132").Value = 1706.5333
$ws.Range("K132").Value = 3158.793
$ws.Range("L132").Value = 5119.5999
$ws.Range("M132").Value = -628.7930000000001
$ws.Range("N132").Value = -10179.5999

# Row 136
$ws.Range("H136").Value = 19610062
$ws.Range("I136").Value = 23810424
$ws.Range("J136").Value = 8378.888999999999
$ws.Range("K136").Value = 71431272
$ws.Range("L136").Value = 25136.667
$ws.Range("M136").Value = -71428722
$ws.Range("N136").Value = -30236.667

Write-Host "Updated 227 cells across 41 rows."
